# Reorder the n95_* and vent_* columns so the *_rank columns are clustered
# together before the value columns.
#
# Before order (per group): percapita, percapita_rank, bytotalcases,
#   bytotalcases_rank, bycases100k, bycases100k_rank
# After order:  percapita_rank, bytotalcases_rank, bycases100k_rank,
#   percapita, bycases100k, bytotalcases   (n95 group, columns G-L)
#   percapita_rank, bytotalcases_rank, bycases100k_rank,
#   percapita, bytotalcases, bycases100k   (vent group, columns N-S)
#
# Column M (ventilators_DELIVERED) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 55

# --- snapshot every source column first (header + all data rows) ---
$g_src = $ws.Range("G1:G$lastRow").Value2
$h_src = $ws.Range("H1:H$lastRow").Value2
$i_src = $ws.Range("I1:I$lastRow").Value2
$j_src = $ws.Range("J1:J$lastRow").Value2
$k_src = $ws.Range("K1:K$lastRow").Value2
$l_src = $ws.Range("L1:L$lastRow").Value2

$n_src = $ws.Range("N1:N$lastRow").Value2
$o_src = $ws.Range("O1:O$lastRow").Value2
$p_src = $ws.Range("P1:P$lastRow").Value2
$q_src = $ws.Range("Q1:Q$lastRow").Value2
$r_src = $ws.Range("R1:R$lastRow").Value2
$s_src = $ws.Range("S1:S$lastRow").Value2

# --- write them back out to their new homes ---
# n95 group: G<-H, H<-J, I<-L, J<-G, K<-K, L<-I
$ws.Range("G1:G$lastRow").Value2 = $h_src
$ws.Range("H1:H$lastRow").Value2 = $j_src
$ws.Range("I1:I$lastRow").Value2 = $l_src
$ws.Range("J1:J$lastRow").Value2 = $g_src
$ws.Range("K1:K$lastRow").Value2 = $k_src
$ws.Range("L1:L$lastRow").Value2 = $i_src

# vent group: N<-O, O<-Q, P<-S, Q<-N, R<-P, S<-R
$ws.Range("N1:N$lastRow").Value2 = $o_src
$ws.Range("O1:O$lastRow").Value2 = $q_src
$ws.Range("P1:P$lastRow").Value2 = $s_src
$ws.Range("Q1:Q$lastRow").Value2 = $n_src
$ws.Range("R1:R$lastRow").Value2 = $p_src
$ws.Range("S1:S$lastRow").Value2 = $r_src
